# Update Linea 141 schedule data (horarios-141) across all three sheets
# generated edit applying the scrape refresh at 12:47:00
$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item('LP1912')
$ws.Range("A2").Value = 'Última actualización: 12:47:00'
$ws.Range("A3").Value = 'Total filas: 189'
$ws.Range("A28").Value = '06:56:24'
$ws.Range("C28").Value = '16_SANTA ANA'
$ws.Range("D28").Value = 25
$ws.Range("A29").Value = '07:15:48'
$ws.Range("C29").Value = '23_HERNANDEZ'
$ws.Range("D29").Value = 6
$ws.Range("C49").Value = '15_ABASTO'
$ws.Range("C50").Value = '11_ETCHEVERRY'
$ws.Range("C53").Value = '10_OLMOS'
$ws.Range("C54").Value = '16_P MOR-SANTA ANA'
$ws.Range("C64").Value = '23_HERNANDEZ'
$ws.Range("C65").Value = '215B_EL PATO'
$ws.Range("A75").Value = '08:40:59'
$ws.Range("C75").Value = '15X38_ABASTO'
$ws.Range("D75").Value = 37
$ws.Range("A76").Value = '08:52:33'
$ws.Range("C76").Value = '14_ABASTO'
$ws.Range("D76").Value = 25
$ws.Range("A77").Value = '08:30:14'
$ws.Range("C77").Value = '27_EL RETIRO'
$ws.Range("D77").Value = 47
$ws.Range("A85").Value = '08:30:14'
$ws.Range("C85").Value = '11_ETCHEVERRY'
$ws.Range("D85").Value = 72
$ws.Range("A86").Value = '08:40:59'
$ws.Range("C86").Value = '16_P MOR-SANTA ANA'
$ws.Range("D86").Value = 62
$ws.Range("C102").Value = '15_ABASTO'
$ws.Range("C103").Value = '14_ABASTO'
$ws.Range("C118").Value = '15X38_ABASTO'
$ws.Range("C119").Value = '14_ABASTO'
$ws.Range("A132").Value = '10:07:51'
$ws.Range("C132").Value = '225_GOMEZ'
$ws.Range("D132").Value = 105
$ws.Range("A133").Value = '11:48:20'
$ws.Range("C133").Value = '23_HERNANDEZ'
$ws.Range("D133").Value = 4
$ws.Range("A147").Value = '11:35:40'
$ws.Range("C147").Value = '23_HERNANDEZ'
$ws.Range("D147").Value = 56
$ws.Range("A148").Value = '11:13:01'
$ws.Range("C148").Value = '27_EL RETIRO'
$ws.Range("D148").Value = 78
$ws.Range("A151").Value = '12:33:54'
$ws.Range("C151").Value = '14_ABASTO'
$ws.Range("D151").Value = 0
$ws.Range("A152").Value = '11:13:01'
$ws.Range("C152").Value = '23_HERNANDEZ'
$ws.Range("D152").Value = 80
$ws.Range("C155").Value = '27_EL RETIRO'
$ws.Range("C156").Value = '23_HERNANDEZ'
$ws.Range("C157").Value = '27_EL RETIRO'
$ws.Range("C158").Value = '23_HERNANDEZ'
$ws.Range("A159").Value = '12:47:00'
$ws.Range("D159").Value = 0
$ws.Range("A161").Value = '12:47:00'
$ws.Range("D161").Value = 1
$ws.Range("A164").Value = '12:47:00'
$ws.Range("D164").Value = 16
$ws.Range("A166").Value = '12:47:00'
$ws.Range("D166").Value = 17
$ws.Range("A168").Value = '12:47:00'
$ws.Range("D168").Value = 30
$ws.Range("A169").Value = '12:47:00'
$ws.Range("B169").Value = '13:19'
$ws.Range("C169").Value = '15_ABASTO'
$ws.Range("D169").Value = 32
$ws.Range("A170").Value = '12:33:54'
$ws.Range("B170").Value = '13:21'
$ws.Range("D170").Value = 48
$ws.Range("A171").Value = '12:47:00'
$ws.Range("B171").Value = '13:22'
$ws.Range("C171").Value = '23_HERNANDEZ'
$ws.Range("D171").Value = 35
$ws.Range("A172").Value = '12:12:04'
$ws.Range("B172").Value = '13:24'
$ws.Range("C172").Value = '23_HERNANDEZ'
$ws.Range("D172").Value = 72
$ws.Range("A173").Value = '12:47:00'
$ws.Range("B173").Value = '13:25'
$ws.Range("C173").Value = '16_P MOR-SANTA ANA'
$ws.Range("D173").Value = 38
$ws.Range("A174").Value = '12:33:54'
$ws.Range("B174").Value = '13:32'
$ws.Range("D174").Value = 59
$ws.Range("A175").Value = '12:33:54'
$ws.Range("B175").Value = '13:32'
$ws.Range("D175").Value = 59
$ws.Range("A176").Value = '12:47:00'
$ws.Range("B176").Value = '13:33'
$ws.Range("C176").Value = '215A_EL PATO'
$ws.Range("D176").Value = 46
$ws.Range("A177").Value = '12:47:00'
$ws.Range("B177").Value = '13:33'
$ws.Range("C177").Value = '14_ABASTO'
$ws.Range("D177").Value = 46
$ws.Range("A178").Value = '12:33:54'
$ws.Range("B178").Value = '13:46'
$ws.Range("C178").Value = '225_GOMEZ'
$ws.Range("D178").Value = 73
$ws.Range("A179").Value = '12:47:00'
$ws.Range("B179").Value = '13:47'
$ws.Range("C179").Value = '225_GOMEZ'
$ws.Range("D179").Value = 60
$ws.Range("A180").Value = '11:55:01'
$ws.Range("B180").Value = '13:49'
$ws.Range("C180").Value = '11_ETCHEVERRY'
$ws.Range("D180").Value = 114
$ws.Range("A181").Value = '12:47:00'
$ws.Range("B181").Value = '13:54'
$ws.Range("C181").Value = '15_ABASTO'
$ws.Range("D181").Value = 67
$ws.Range("A182").Value = '12:47:00'
$ws.Range("B182").Value = '13:58'
$ws.Range("C182").Value = '16_SANTA ANA'
$ws.Range("D182").Value = 71
$ws.Range("B183").Value = '14:01'
$ws.Range("C183").Value = '10_OLMOS'
$ws.Range("D183").Value = 88
$ws.Range("A184").Value = '12:47:00'
$ws.Range("B184").Value = '14:02'
$ws.Range("C184").Value = '10_OLMOS'
$ws.Range("D184").Value = 75
$ws.Range("E184").Value = 'LP1912'
$ws.Range("A185").Value = '12:47:00'
$ws.Range("B185").Value = '14:07'
$ws.Range("C185").Value = '23_HERNANDEZ'
$ws.Range("D185").Value = 80
$ws.Range("E185").Value = 'LP1912'
$ws.Range("A186").Value = '12:33:54'
$ws.Range("B186").Value = '14:16'
$ws.Range("C186").Value = '27_EL RETIRO'
$ws.Range("D186").Value = 103
$ws.Range("E186").Value = 'LP1912'
$ws.Range("A187").Value = '12:33:54'
$ws.Range("B187").Value = '14:17'
$ws.Range("C187").Value = '11_ETCHEVERRY'
$ws.Range("D187").Value = 104
$ws.Range("E187").Value = 'LP1912'
$ws.Range("A188").Value = '12:47:00'
$ws.Range("B188").Value = '14:17'
$ws.Range("C188").Value = '27_EL RETIRO'
$ws.Range("D188").Value = 90
$ws.Range("E188").Value = 'LP1912'
$ws.Range("A189").Value = '12:47:00'
$ws.Range("B189").Value = '14:18'
$ws.Range("C189").Value = '11_ETCHEVERRY'
$ws.Range("D189").Value = 91
$ws.Range("E189").Value = 'LP1912'
$ws.Range("A190").Value = '12:47:00'
$ws.Range("B190").Value = '14:27'
$ws.Range("C190").Value = '16_SANTA ANA'
$ws.Range("D190").Value = 100
$ws.Range("E190").Value = 'LP1912'
$ws.Range("A191").Value = '12:33:54'
$ws.Range("B191").Value = '14:31'
$ws.Range("C191").Value = '14X44_ABASTO'
$ws.Range("D191").Value = 118
$ws.Range("E191").Value = 'LP1912'
$ws.Range("A192").Value = '12:47:00'
$ws.Range("B192").Value = '14:32'
$ws.Range("C192").Value = '14X44_ABASTO'
$ws.Range("D192").Value = 105
$ws.Range("E192").Value = 'LP1912'
$ws.Range("A193").Value = '12:47:00'
$ws.Range("B193").Value = '14:34'
$ws.Range("C193").Value = '215C_EL PATO'
$ws.Range("D193").Value = 107
$ws.Range("E193").Value = 'LP1912'
$ws.Range("A194").Value = '12:47:00'
$ws.Range("B194").Value = '14:39'
$ws.Range("C194").Value = '16_P MOR-SANTA ANA'
$ws.Range("D194").Value = 112
$ws.Range("E194").Value = 'LP1912'

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Range("A2").Value = 'Última actualización: 12:47:00'
$ws.Range("A3").Value = 'Total filas: 26'
$ws.Range("A28").Value = '12:47:00'
$ws.Range("D28").Value = 17
$ws.Range("A30").Value = '12:47:00'
$ws.Range("D30").Value = 46
$ws.Range("A31").Value = '12:47:00'
$ws.Range("B31").Value = '14:34'
$ws.Range("C31").Value = '215C_EL PATO'
$ws.Range("D31").Value = 107
$ws.Range("E31").Value = 'LP1912'

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Range("A2").Value = 'Última actualización: 12:47:00'
$ws.Range("A19").Value = '08:52:33'
$ws.Range("C19").Value = '215B_LP-P MOR-1 Y 57'
$ws.Range("D19").Value = 98
$ws.Range("A20").Value = '10:07:51'
$ws.Range("C20").Value = '215A_LA PLATA'
$ws.Range("D20").Value = 23
$ws.Range("A25").Value = '12:47:00'
$ws.Range("D25").Value = 25
$ws.Range("A27").Value = '12:47:00'
$ws.Range("D27").Value = 34
$ws.Range("A29").Value = '12:47:00'
$ws.Range("D29").Value = 70

